$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the "Fecha" (D) values between the first pair of rows (2,3) and the
# second pair of rows (6,7): rows 2 & 3 move from 44559 -> 44574, while
# rows 6 & 7 move from 44574 -> 44559.
$ws.Range("D2").Value = 44574
$ws.Range("D3").Value = 44574
$ws.Range("D6").Value = 44559
$ws.Range("D7").Value = 44559
